$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = ". . . es geschehen viele sehr wichtige Dinge in der Welt, über die die Öffentlichkeit nie informiert wird."
$ws.Range("B3").Value = ". . . Politiker geben uns normalerweise keine Auskunft über die wahren Motive ihrer Entscheidungen."
$ws.Range("B4").Value = ". . . Regierungsbehörden überwachen alle Bürger genau."
$ws.Range("B5").Value = ". . . Ereignisse, die auf den ersten Blick nicht miteinander in Verbindung zu stehen scheinen, sind oft das Ergebnis geheimer Aktivitäten."
$ws.Range("B6").Value = ". . . es gibt geheime Organisationen, die großen Einfluss auf politische Entscheidungen haben."
